$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are textual (may contain multiple "." as thousands
# separators); force Text number format first so Excel does not coerce them
# into numeric values and strip formatting (e.g. trailing zeros).
$dCells = @("D2","D3","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D18","D21","D22","D23","D24","D26","D27","D28","D29","D30","D31","D32","D35","D37","D38","D39","D40","D41","D43","D46","D47","D49")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.374.89'
$ws.Range("D3").Value = '1.840.60'
$ws.Range("D5").Value = '315.16'
$ws.Range("D6").Value = '1.013'
$ws.Range("D7").Value = '0.4745'
$ws.Range("D8").Value = '0.3702'
$ws.Range("D9").Value = '0.07468'
$ws.Range("D10").Value = '0.8855'
$ws.Range("D11").Value = '20.51'
$ws.Range("D12").Value = '1.845.91'
$ws.Range("D13").Value = '0.07386'
$ws.Range("D14").Value = '5.487'
$ws.Range("D15").Value = '93.26'
$ws.Range("D16").Value = '6.580'
$ws.Range("D18").Value = '0.000008851'
$ws.Range("D21").Value = '27.412.72'
$ws.Range("D22").Value = '5.359'
$ws.Range("D23").Value = '10.73'
$ws.Range("D24").Value = '2.065.25'
$ws.Range("D26").Value = '152.08'
$ws.Range("D27").Value = '18.66'
$ws.Range("D28").Value = '2.184'
$ws.Range("D29").Value = '5.283'
$ws.Range("D30").Value = '118.01'
$ws.Range("D31").Value = '0.08966'
$ws.Range("D32").Value = '0.7622'
$ws.Range("D35").Value = '2.941'
$ws.Range("D37").Value = '1.106'
$ws.Range("D38").Value = '0.05371'
$ws.Range("D39").Value = '0.01965'
$ws.Range("D40").Value = '3.004'
$ws.Range("D41").Value = '7.307'
$ws.Range("D43").Value = '2.379'
$ws.Range("D46").Value = '0.4991'
$ws.Range("D47").Value = '10.52'
$ws.Range("D49").Value = '105.23'

# Volume(1h) (column E) values
$ws.Range("E2").Value = '  +1.88%  '
$ws.Range("E3").Value = '  +1.45%  '
$ws.Range("E4").Value = '  +1.37%  '
$ws.Range("E5").Value = '  +1.99%  '
$ws.Range("E6").Value = '  +1.16%  '
$ws.Range("E7").Value = '  +1.85%  '
$ws.Range("E9").Value = '  +1.41%  '
$ws.Range("E10").Value = '  +1.77%  '
$ws.Range("E11").Value = '  +0.62%  '
$ws.Range("E12").Value = '  +1.62%  '
$ws.Range("E13").Value = '  +4.56%  '
$ws.Range("E14").Value = '  +2.27%  '
$ws.Range("E15").Value = '  +1.92%  '
$ws.Range("E16").Value = '  +1.16%  '
$ws.Range("E17").Value = '  +1.19%  '
$ws.Range("E18").Value = '  +1.90%  '
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("E20").Value = '  +0.97%  '
$ws.Range("E21").Value = '  +1.93%  '
$ws.Range("E22").Value = '  +0.43%  '
$ws.Range("E23").Value = '  +1.45%  '
$ws.Range("E24").Value = '  +0.99%  '
$ws.Range("E25").Value = '  +0.47%  '
$ws.Range("E26").Value = '  +1.13%  '
$ws.Range("E27").Value = '  +1.91%  '
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("E29").Value = '  -0.67%  '
$ws.Range("E30").Value = '  +2.01%  '
$ws.Range("E31").Value = '  +0.45%  '
$ws.Range("E32").Value = '  -0.49%  '
$ws.Range("E33").Value = '  +1.37%  '
$ws.Range("E34").Value = '  +1.33%  '
$ws.Range("E35").Value = '  +1.32%  '
$ws.Range("E36").Value = '  +1.24%  '
$ws.Range("E37").Value = '  +1.77%  '
$ws.Range("E38").Value = '  +1.72%  '
$ws.Range("E39").Value = '  +0.23%  '
$ws.Range("E40").Value = '  +2.55%  '
$ws.Range("E41").Value = '  +0.82%  '
$ws.Range("E42").Value = '  +0.83%  '
$ws.Range("E43").Value = '  +1.54%  '
$ws.Range("E44").Value = '  +0.51%  '
$ws.Range("E45").Value = '  +1.56%  '
$ws.Range("E46").Value = '  +1.38%  '
$ws.Range("E47").Value = '  +0.67%  '
$ws.Range("E49").Value = '  +1.51%  '
$ws.Range("E50").Value = '  +0.52%  '
$ws.Range("E51").Value = '  +0.76%  '
